$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the company name (razon social)
$ws.Range("E7").Value = "CENTRO DE DIAGNÓSTICO AUTOMOTOR REVICAR LTDA"

# Update total "Valor Mora"
$ws.Range("E11").Value = 72682

# Update worker/period counts
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 2

# Remove the first worker's row (WILSON JESUS ACENDRA RODRIGUEZ / period 1803).
# This shifts the remaining two TANIA JULIO SARMIENTO rows (periods 2107, 2106)
# up into rows 16 and 17, and shifts the signature-block rows up to 22/23.
$ws.Rows.Item(16).Delete()

# Update the remaining worker rows' amounts
$ws.Cells.Item(16, 6).Value = 36341
$ws.Cells.Item(16, 7).Value = 908526
$ws.Cells.Item(17, 6).Value = 36341
$ws.Cells.Item(17, 7).Value = 908526

# After the deletion, row 16 holds period 2107 and row 17 holds period 2106;
# the target order is 2106 first, then 2107 - swap the period values.
$period16 = $ws.Cells.Item(16, 5).Value()
$period17 = $ws.Cells.Item(17, 5).Value()
$ws.Cells.Item(16, 5).Value = $period17
$ws.Cells.Item(17, 5).Value = $period16
